$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1) ---
$ws.Range("A1").Value = "Fecha"
$ws.Range("B1").Value = "Hora"
$ws.Range("C1").Value = "WC47 NACP"
$ws.Range("D1").Value = "WC48 P5F"
$ws.Range("E1").Value = "WC49 P5H"
$ws.Range("F1").Value = "WV50 FILTER"
$ws.Range("G1").Value = "SPL"

# --- Append new data rows 31-69 ---
$newData = @(
    @("2024-05-13","11:59:57","-","-","Atasco tuerca","-","-"),
    @("2024-05-13","12:05:10","Fallo tornillo","-","-","-","-"),
    @("2024-05-13","12:05:14","-","-","Atasco tuerca","-","-"),
    @("2024-05-13","12:05:18","-","-","-","Robot no coloca bien ferrita","-"),
    @("2024-05-13","12:06:45","-","Cámara no detecta Pcb","-","-","-"),
    @("2024-05-13","12:13:21","-","-","-","Robot no coge PCB","-"),
    @("2024-05-13","12:13:27","No coge placa","-","-","-","-"),
    @("2024-05-13","12:22:33","-","-","-","Robot no coge PCB","-"),
    @("2024-05-13","12:22:36","-","-","-","-","Marco atascado en parte inferior"),
    @("2024-05-13","12:25:17","-","Fallo etiqueta","-","-","-"),
    @("2024-05-13","12:25:27","No atornilla clips","-","-","-","-"),
    @("2024-05-13","12:25:57","-","No detecta presencia power CP","-","-","-"),
    @("2024-05-13","12:26:01","-","-","-","No coloca bien la pcb","-"),
    @("2024-05-13","12:26:05","-","-","-","-","Error en sensor de salida"),
    @("2024-05-13","12:26:24","Fallo tornillo","-","-","-","-"),
    @("2024-05-13","12:26:27","-","Cámara no detecta Power CP","-","-","-"),
    @("2024-05-13","12:26:30","-","-","No coloca bien el sealling","-","-"),
    @("2024-05-13","12:26:34","-","-","-","Fallo visión core","-"),
    @("2024-05-13","12:29:16","No coge placa","-","-","-","-"),
    @("2024-05-13","12:29:19","No atornilla clips","-","-","-","-"),
    @("2024-05-13","12:29:22","-","Cámara no detecta Top cover","-","-","-"),
    @("2024-05-13","12:29:27","-","-","No lee QR","-","-"),
    @("2024-05-13","12:29:32","-","-","-","-","Error en sensor de salida"),
    @("2024-05-13","12:29:39","No atornilla clips","-","-","-","-"),
    @("2024-05-13","12:29:41","No atornilla clips","-","-","-","-"),
    @("2024-05-13","12:29:43","No atornilla clips","-","-","-","-"),
    @("2024-05-13","12:32:27","Secuencia atornillador","-","-","-","-"),
    @("2024-05-13","12:32:32","Fallo fijador tapa","-","-","-","-"),
    @("2024-05-13","12:32:34","Fallo fijador tapa","-","-","-","-"),
    @("2024-05-13","12:32:38","-","-","-","-","Marco atascado en parte inferior"),
    @("2024-05-13","12:36:50","Fallo tolva","-","-","-","-"),
    @("2024-05-13","12:36:56","Palet atascado en la curva","-","-","-","-"),
    @("2024-05-13","12:37:05","-","-","Etiquetadora","-","-"),
    @("2024-05-13","12:37:09","-","-","-","Core enganchado","-"),
    @("2024-05-13","12:37:21","Secuencia atornillador","-","-","-","-"),
    @("2024-05-13","12:37:24","-","-","-","Robot no coloca bien ferrita","-"),
    @("2024-05-13","12:37:37","-","-","-","Fallo atornillador","-"),
    @("2024-05-13","12:38:08","-","No detecta presencia power CP","-","-","-"),
    @("2024-05-13","12:38:17","-","Cámara no detecta busbar","-","-","-")
)

$startRow = 31
for ($i = 0; $i -lt $newData.Count; $i++) {
    $rowIndex = $startRow + $i
    $rowValues = $newData[$i]
    for ($col = 0; $col -lt $rowValues.Count; $col++) {
        $cell = $ws.Cells.Item($rowIndex, $col + 1)
        $value = $rowValues[$col]
        if ($col -eq 0) {
            # Column A holds date-like text (e.g. 2024-05-13) that Excel would
            # otherwise auto-convert to a date serial number. Force plain text
            # by using a leading quote prefix, then reset the cell style so no
            # residual formatting/quote-prefix marker is left behind.
            $cell.Value = "'" + $value
            $cell.Style = "Normal"
        } else {
            $cell.Value = $value
        }
    }
}
